$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 149, shifting existing rows 149:202 down to 150:203
$ws.Rows.Item(149).Insert()

# Populate the new row 149 with the new record
$ws.Range("A149").Value = 3
$ws.Range("B149").Value = "Femacal de La Calera"
$ws.Range("C149").Value = "Coquimbo"
$ws.Range("D149").Value = 44468
$ws.Range("E149").Value = 5
$ws.Range("F149").Value = 100112043
$ws.Range("G149").Value = "Pepino ensalada"
$ws.Range("H149").Value = "Sin especificar"
$ws.Range("I149").Value = "Primera"
$ws.Range("J149").Value = 130
$ws.Range("K149").Value = 13000
$ws.Range("L149").Value = 14000
$ws.Range("M149").Value = 13538
$ws.Range("N149").Value = "$/caja 70 unidades"
$ws.Range("O149").Value = "Región de Arica y Parinacota"
$ws.Range("P149").Value = 193
$ws.Range("Q149").Value = 70
$ws.Range("R149").Value = "Hortaliza"

# Apply the same date number format (style index 2 -> numFmtId 165) used by other D-column cells
$ws.Range("D149").NumberFormat = $ws.Range("D150").NumberFormat

Write-Host "Row inserted and populated"
